# Apply cryptos.xlsx price/volume update (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.640.98"
$ws.Range("E2").Value = "  -2.51%  "

$ws.Range("D3").Value = "2.399.49"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("D4").Value = "'0.993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.66%  "

$ws.Range("D5").Value = "'570.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "'139.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.16%  "

$ws.Range("D9").Value = "2.379.62"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  -2.35%  "

$ws.Range("E13").Value = "  -1.36%  "

$ws.Range("D14").Value = "'25.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "'0.0000169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").Value = "60.639.99"
$ws.Range("E17").Value = "  -2.46%  "

$ws.Range("D18").Value = "2.364.43"
$ws.Range("E18").Value = "  -2.29%  "

$ws.Range("D19").Value = "'10.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.74%  "

$ws.Range("D20").Value = "'7.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("D21").Value = "'321.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("D23").Value = "'6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("E25").Value = "  -5.52%  "

$ws.Range("D26").Value = "'64.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("E27").Value = "  -8.19%  "

$ws.Range("D28").Value = "'570.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.99%  "

$ws.Range("D29").Value = "2.486.90"
$ws.Range("E29").Value = "  -3.08%  "

$ws.Range("E30").Value = "  -4.16%  "

$ws.Range("D31").Value = "'7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("D32").Value = "'1.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.63%  "

$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("E34").Value = "  -5.52%  "

$ws.Range("D35").Value = "'1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "'4.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.82%  "

$ws.Range("D37").Value = "'0.367"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").Value = "'147.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.38%  "

$ws.Range("D40").Value = "'18.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").Value = "'5.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.23%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "'1.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.56%  "

$ws.Range("D44").Value = "'40.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.59%  "

$ws.Range("D45").Value = "'2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.39%  "

$ws.Range("E46").Value = "  +19.84%  "

$ws.Range("D47").Value = "'139.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("E48").Value = "  -3.50%  "

$ws.Range("D49").Value = "'0.583"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.75%  "

$ws.Range("D50").Value = "'0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.00%  "

$ws.Range("D51").Value = "'19.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
